# Update with RQ1 results: the "6+" authors bucket is now "11+".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 previously held the "6+" label in column A; relabel to "11+".
$ws.Range("A4").Value = "11+"
$ws.Range("A5").Value = "11+"

# Leave the final selection on A5, matching the saved sheet view.
$ws.Range("A5").Select()
